$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix total marks: "Marking" row (per-question marks) corrections
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Fix total marks: "Total" row corrections
$ws.Range("B12").Value = 48
$ws.Range("C12").Value = -6

# Fix displayed "obtained / max" summary text
$ws.Range("E12").Value = "42 / 112"
